$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert three new rows after "Device Mode" (row 4): Temp Correction,
#    Temp Interval, Relay state. Each insert shifts everything below it down
#    by one row; we copy cell formatting (style only) from row 4 so the new
#    rows pick up the existing bordered style (s="2") instead of the default.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).EntireRow.Insert()
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

$ws.Rows.Item(6).EntireRow.Insert()
$ws.Range("A4:D4").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)

$ws.Rows.Item(7).EntireRow.Insert()
$ws.Range("A4:D4").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Insert one new row for "MQTT Topic" right after "MQTT Password" (which,
#    post-shift, now lives at row 14). Copy formatting from row 14.
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).EntireRow.Insert()
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# The row-insert above duplicates row 14's formatting across the whole row,
# including the styled (but otherwise empty) J14 cell, giving a stray J15.
# Only J14 should carry that style, so remove the J15 cell entirely.
$ws.Range("J15").Clear()

# ---------------------------------------------------------------------------
# 3) Remove one filler row from the bottom: four rows were inserted above but
#    only three net rows are new table entries (one old filler row is no
#    longer needed), bringing the sheet from 28 to 31 rows total.
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).EntireRow.Delete()

# ---------------------------------------------------------------------------
# 4) Write the final field values/labels for rows 4-16 and make sure every
#    Start/End formula explicitly chains off the row immediately above it
#    (row-insert preserves old absolute references rather than re-deriving
#    the "row above" relationship, so these must be set explicitly).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Device Mode"
$ws.Range("B4").Value = 1
$ws.Range("C4").Formula = "=D3+1"
$ws.Range("D4").Formula = "=B4+C4-1"

$ws.Range("A5").Value = "Temp Correction"
$ws.Range("B5").Value = 5
$ws.Range("C5").Formula = "=D4+1"
$ws.Range("D5").Formula = "=B5+C5-1"

$ws.Range("A6").Value = "Temp Interval"
$ws.Range("B6").Value = 8
$ws.Range("C6").Formula = "=D5+1"
$ws.Range("D6").Formula = "=B6+C6-1"

$ws.Range("A7").Value = "Relay state"
$ws.Range("B7").Value = 1
$ws.Range("C7").Formula = "=D6+1"
$ws.Range("D7").Formula = "=B7+C7-1"

$ws.Range("A8").Value = "Reserved"
$ws.Range("B8").Value = 113
$ws.Range("C8").Formula = "=D7+1"
$ws.Range("D8").Formula = "=B8+C8-1"

$ws.Range("A9").Value = "WIFI SSID"
$ws.Range("B9").Value = 32
$ws.Range("C9").Formula = "=D8+1"
$ws.Range("D9").Formula = "=B9+C9-1"

$ws.Range("A10").Value = "WIFI Password"
$ws.Range("B10").Value = 32
$ws.Range("C10").Formula = "=D9+1"
$ws.Range("D10").Formula = "=B10+C10-1"

$ws.Range("A11").Value = "MQTT Host"
$ws.Range("B11").Value = 32
$ws.Range("C11").Formula = "=D10+1"
$ws.Range("D11").Formula = "=B11+C11-1"

$ws.Range("A12").Value = "MQTT Port"
$ws.Range("B12").Value = 5
$ws.Range("C12").Formula = "=D11+1"
$ws.Range("D12").Formula = "=B12+C12-1"

$ws.Range("A13").Value = "MQTT User"
$ws.Range("B13").Value = 32
$ws.Range("C13").Formula = "=D12+1"
$ws.Range("D13").Formula = "=B13+C13-1"

$ws.Range("A14").Value = "MQTT Password"
$ws.Range("B14").Value = 32
$ws.Range("C14").Formula = "=D13+1"
$ws.Range("D14").Formula = "=B14+C14-1"

$ws.Range("A15").Value = "MQTT Topic"
$ws.Range("B15").Value = 32
$ws.Range("C15").Formula = "=D14+1"
$ws.Range("D15").Formula = "=B15+C15-1"

$ws.Range("A16").Value = "Reserved"
$ws.Range("B16").Value = 32
$ws.Range("C16").Formula = "=D15+1"
$ws.Range("D16").Formula = "=B16+C16-1"

# ---------------------------------------------------------------------------
# 5) Move the active-cell selection to A7, matching the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("A7").Select() | Out-Null
